$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row -------------------------------------------
# "<Name>_old" -> "<Name>_FV2410"  and  "<Name>_new" -> "<Name>_FV2504"
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value  = $baseNames[$i] + "_FV2410"
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2504"
}
# column K ("diff") keeps its name

# --- 2) Turn the data range into an Excel Table --------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row --------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
